$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 2 (B2:G2)
$ws.Range("B2").Value = -0.1338314788754218
$ws.Range("C2").Value = 2.04803740314168
$ws.Range("D2").Value = 20.30227337679758
$ws.Range("E2").Value = 4.505804409514197
$ws.Range("F2").Value = 4.605038435222934
$ws.Range("G2").Value = 23

# Row 3 now takes old row 2's values
$ws.Range("B3").Value = -0.2196048395615291
$ws.Range("C3").Value = 2.254442816470436
$ws.Range("D3").Value = 18.45943047502218
$ws.Range("E3").Value = 4.296443933652827
$ws.Range("F3").Value = 4.39180238609223
$ws.Range("G3").Value = 22

# Row 4 now takes old row 3's values
$ws.Range("B4").Value = -0.6722431955569657
$ws.Range("C4").Value = 1.734626294717863
$ws.Range("D4").Value = 10.27176916279258
$ws.Range("E4").Value = 3.204960087550636
$ws.Range("F4").Value = 3.211051410560245
$ws.Range("G4").Value = 21

# Row 5 now takes old row 4's values
$ws.Range("B5").Value = -0.2619961853625158
$ws.Range("C5").Value = 1.574840796137216
$ws.Range("D5").Value = 11.52157846430296
$ws.Range("E5").Value = 3.394345071483298
$ws.Range("F5").Value = 3.472135163382664
$ws.Range("G5").Value = 20

# Row 6 now takes old row 5's values
$ws.Range("B6").Value = -0.2304210235372018
$ws.Range("C6").Value = 1.666997890576107
$ws.Range("D6").Value = 10.06653556874298
$ws.Range("E6").Value = 3.172780416092955
$ws.Range("F6").Value = 3.251114276439574
$ws.Range("G6").Value = 19

# Row 7 now takes old row 6's values
$ws.Range("B7").Value = -0.2737229439002919
$ws.Range("C7").Value = 1.775059182383401
$ws.Range("D7").Value = 10.57130412137001
$ws.Range("E7").Value = 3.251354198079626
$ws.Range("F7").Value = 3.333738739228395
$ws.Range("G7").Value = 18

# Row 8 now takes old row 7's values
$ws.Range("B8").Value = -0.1641233777288165
$ws.Range("C8").Value = 1.781553791812374
$ws.Range("D8").Value = 11.51292075854603
$ws.Range("E8").Value = 3.393069518672736
$ws.Range("F8").Value = 3.493402108638973
$ws.Range("G8").Value = 17

# Row 9 now takes old row 8's values
$ws.Range("B9").Value = -0.1510345969195566
$ws.Range("C9").Value = 1.85455728114967
$ws.Range("D9").Value = 11.36885938571983
$ws.Range("E9").Value = 3.37177392268815
$ws.Range("F9").Value = 3.478857733797603
$ws.Range("G9").Value = 16

# Row 10 now takes old row 9's values
$ws.Range("B10").Value = -0.1166201009408896
$ws.Range("C10").Value = 1.99642669408382
$ws.Range("D10").Value = 12.99001923899349
$ws.Range("E10").Value = 3.604166927182131
$ws.Range("F10").Value = 3.728713727525796
$ws.Range("G10").Value = 15

# Row 11 now takes old row 10's values
$ws.Range("B11").Value = -0.07784734317373232
$ws.Range("C11").Value = 1.990677678875613
$ws.Range("D11").Value = 13.16673245956295
$ws.Range("E11").Value = 3.628599242071649
$ws.Range("F11").Value = 3.764708707805368
$ws.Range("G11").Value = 14
